$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.098.50'
$ws.Range('E2').Value = '  -1.17%  '

# Row 3
$ws.Range('D3').Value = '1.791.71'
$ws.Range('E3').Value = '  -0.26%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '225.04'
$ws.Range('E5').Value = '  +0.90%  '

# Row 6
$ws.Range('E6').Value = '  -0.26%  '

# Row 7
$ws.Range('E7').Value = '  +0.06%  '

# Row 8
$ws.Range('E8').Value = '  +0.97%  '

# Row 9
$ws.Range('E9').Value = '  -1.13%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0707'
$ws.Range('E10').Value = '  +0.25%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').Value = '  +0.08%  '

# Row 12
$ws.Range('D12').Value = '2.050.26'
$ws.Range('E12').Value = '  -0.17%  '

# Row 13
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.797.96'
$ws.Range('E13').Value = '  -0.05%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.82'
$ws.Range('E14').Value = '  -1.26%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.624'
$ws.Range('E15').Value = '  -2.21%  '

# Row 16
$ws.Range('D16').Value = '34.055.74'
$ws.Range('E16').Value = '  -1.30%  '

# Row 17
$ws.Range('E17').Value = '  -2.36%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.06'
$ws.Range('E18').Value = '  -1.46%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.57'
$ws.Range('E19').Value = '  -2.99%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0786'
$ws.Range('E20').Value = '  -1.52%  '

# Row 21
$ws.Range('E21').Value = '  +0.02%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.68'
$ws.Range('E22').Value = '  -3.50%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.10'
$ws.Range('E23').Value = '  -2.98%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  -2.80%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.31'
$ws.Range('E25').Value = '  -1.31%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.26'
$ws.Range('E26').Value = '  -0.58%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.04'
$ws.Range('E27').Value = '  -1.20%  '

# Row 28
$ws.Range('E28').Value = '  -1.49%  '

# Row 29
$ws.Range('E29').Value = '  +0.09%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0518'
$ws.Range('E30').Value = '  -1.00%  '

# Row 31
$ws.Range('E31').Value = '  +2.29%  '

# Row 32
$ws.Range('E32').Value = '  -2.69%  '

# Row 33
$ws.Range('E33').Value = '  -2.15%  '

# Row 34
$ws.Range('E34').Value = '  -3.41%  '

# Row 35
$ws.Range('D35').Value = '1.395.08'
$ws.Range('E35').Value = '  -1.83%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.646'
$ws.Range('E36').Value = '  +1.77%  '

# Row 37
$ws.Range('E37').Value = '  -1.19%  '

# Row 38
$ws.Range('E38').Value = '  -2.45%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.20'
$ws.Range('E39').Value = '  +3.61%  '

# Row 40
$ws.Range('B40').Value = 'HuobiToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.35'
$ws.Range('E40').Value = '  -0.16%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '78.99'
$ws.Range('E41').Value = '  -4.44%  '

# Row 42
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.70'
$ws.Range('E42').Value = '  -2.97%  '

# Row 43
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.915'
$ws.Range('E43').Value = '  -3.84%  '

# Row 44
$ws.Range('D44').Value = '0.0₆0146'
$ws.Range('E44').Value = '  +17.92%  '

# Row 45
$ws.Range('E45').Value = '  +0.66%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.25'
$ws.Range('E46').Value = '  +2.69%  '

# Row 47
$ws.Range('E47').Value = '  -0.50%  '

# Row 48
$ws.Range('E48').Value = '  -2.27%  '

# Row 49
$ws.Range('D49').Value = '1.949.50'
$ws.Range('E49').Value = '  +0.13%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '12.08'
$ws.Range('E50').Value = '  -1.60%  '

# Row 51
$ws.Range('E51').Value = '  -0.03%  '
